$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80:B80").Copy()
$ws.Range("A81:B81").PasteSpecial(-4122)

$ws.Range("A81").Value = "UserNotActive"
$ws.Range("B81").Value = "Kullanıcı hesabınız aktif değildir. Lütfen TestOkur yetkilileri ile görüşünüz"
